# Add "audios" clean up prompt rows: every trigger row gets an extra
# "p_<EVENT>" / ="<EVENT>"; // <audio prompt text> style expression,
# and several of the audio prompts are split into red/blue partner
# specific variants.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 2-13. Column layout:
#   A = p_<EVENT>            (new prompt identifier)
#   B = "="                  (unchanged literal)
#   C = '"'                  (unchanged literal)
#   D = <EVENT>               (trigger / event name)
#   E = '"'                  (unchanged literal)
#   F = ";"                  (was "," before the edit)
#   G = "//"                 (unchanged literal)
#   H = <prompt text>         (audio / on-screen prompt)
$rows = @(
    @{ Row = 2;  Event = 'ON_VERTICAL_STARTUP';   Prompt = 'Please place the iPad flat in front of you to begin.' },
    @{ Row = 3;  Event = 'ON_STARTUP';             Prompt = 'Welcome to MC Partner. Is this the first time you have played? Tap the green Yes or the red No' },
    @{ Row = 4;  Event = 'ON_INTRO';               Prompt = 'This game will help you learn with your partner. One of you will be the red partner and one will be the blue partner. You decide!' },
    @{ Row = 5;  Event = 'RED_PICK_UP_IPAD';       Prompt = 'Red partner now pick up the iPad and hold it so you can see my face but your partner can''t.' },
    @{ Row = 6;  Event = 'BLUE_PICK_UP_IPAD';      Prompt = 'Blue partner now pick up the iPad and hold it so you can see my face but your partner can''t.' },
    @{ Row = 7;  Event = 'READ_THE_WORD';          Prompt = 'Great! Now read the word you see out loud so your partner can hear it. When you have read the word, place the iPad flat in front of you again.' },
    @{ Row = 8;  Event = 'PUT_DOWN_THE_IPAD';      Prompt = 'When you have read the word, place the iPad flat in front of you again.' },
    @{ Row = 9;  Event = 'RED_SELECT_STIMULUS';    Prompt = 'Red partner now find the word your partner read to you and tap it with your finger.' },
    @{ Row = 10; Event = 'BLUE_SELECT_STIMULUS';   Prompt = 'Blue partner now find the word your partner read to you and tap it with your finger.' },
    @{ Row = 11; Event = 'CORRECT_STIMULUS';       Prompt = 'The right answer was [ANSWER]' },
    @{ Row = 12; Event = 'INCORRECT_SELECTION';    Prompt = 'The word you chose was [ANSWER]' },
    @{ Row = 13; Event = 'CORRECT_SELECTION';      Prompt = 'Great! The right answer was [ANSWER]' }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $ws.Range("A$rowNum").Value = "p_" + $r.Event
    # A leading "=" would otherwise be parsed by Excel as the start of a
    # formula, so force it to be stored as a literal piece of text, then
    # drop the resulting quote-prefix formatting so the cell keeps the
    # default style (no explicit s="..." attribute).
    $ws.Range("B$rowNum").Value = "'="
    $ws.Range("B$rowNum").Style = "Normal"
    $ws.Range("C$rowNum").Value = '"'
    $ws.Range("D$rowNum").Value = $r.Event
    $ws.Range("E$rowNum").Value = '"'
    $ws.Range("F$rowNum").Value = ";"
    $ws.Range("G$rowNum").Value = "//"
    $ws.Range("H$rowNum").Value = $r.Prompt
}

# The sheet used to contain an entirely empty row 14 (between the data
# block ending at row 13 and the two height-only formatting rows that
# were at 16 and 20). Removing it shifts those rows up to 15 and 19,
# matching the updated layout.
$ws.Rows.Item(14).Delete()

# Update the sheet view: drop the stale top-left/selection anchor left
# over from scrolling, and leave the selection on F16 instead.
$ws.Range("F16").Select()
